# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 555
$ws1.Range("G3").Value = 70

$ws1.Range("F4").Value = 1581
$ws1.Range("G4").Value = 70

$ws1.Range("F5").Value = 172
$ws1.Range("G5").Value = 58

$ws1.Range("F7").Value = 5212

$ws1.Range("F8").Value = 194

$ws1.Range("F11").Value = 68

$ws1.Range("F16").Value = 6566

$ws1.Range("F19").Value = 141

$ws1.Range("F22").Value = 1015

$ws1.Range("F23").Value = 15697

$ws1.Range("F25").Value = 20

$ws1.Range("F29").Value = 11152

$ws1.Range("F30").Value = 793

$ws1.Range("F32").Value = 267

$ws1.Range("F34").Value = 25

$ws1.Range("F35").Value = 311

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 555
$ws4.Range("G3").Value = 70

$ws4.Range("F4").Value = 1581
$ws4.Range("G4").Value = 70

$ws4.Range("F5").Value = 172
$ws4.Range("G5").Value = 58

$ws4.Range("F8").Value = 5212

$ws4.Range("F9").Value = 194

$ws4.Range("F13").Value = 68

$ws4.Range("F19").Value = 6566

$ws4.Range("F22").Value = 141

$ws4.Range("F26").Value = 1015

$ws4.Range("F27").Value = 15697

$ws4.Range("F29").Value = 20

$ws4.Range("F34").Value = 11152

$ws4.Range("F35").Value = 793

$ws4.Range("F37").Value = 267

$ws4.Range("F39").Value = 25

$ws4.Range("F40").Value = 311
